$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data appended below the existing table (row 7)
$row = 7

$ws.Cells.Item($row, 1).Value = "DF"            # A7 - UF
$ws.Cells.Item($row, 2).Value = "00X0098"        # B7 - FRU
# C7, D7, E7 (SUB1/SUB2/SUB3) stay blank, matching the other data rows
$ws.Cells.Item($row, 6).Value = "TESTEOLOBATO"   # F7 - DESCRICAO
$ws.Cells.Item($row, 7).Value = "T"              # G7 - MAQUINAS
$ws.Cells.Item($row, 8).Value = "T - (T 01/11/25_12H) - DF"  # H7 - CLIENTE

# I7 (DATA_FIM) must stay as literal text "01/11/25", not get auto-converted
# to a date serial by Excel's smart entry, so force text format first.
$cellI = $ws.Cells.Item($row, 9)
$cellI.NumberFormat = "@"
$cellI.Value = "01/11/25"
$cellI.Style = "Normal"

$ws.Cells.Item($row, 10).Value = "12H"           # J7 - SLA
$ws.Cells.Item($row, 11).Value = "14/11/25"      # K7 - DATA_VERIFICACAO
$ws.Cells.Item($row, 12).Value = "DENTRO"        # L7 - STATUS
# M7 (DATA_FIM_DT) stays blank, same as M6
